# Apply updated odds values to "Jogos da Semana" FlashScore workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 updates ---
$ws.Range("G6").Value  = 2.32
$ws.Range("I6").Value  = 3.2
$ws.Range("M6").Value  = 1.07
$ws.Range("N6").Value  = 7
$ws.Range("O6").Value  = 1.47
$ws.Range("Q6").Value  = 1.84
$ws.Range("R6").Value  = 1.89
$ws.Range("X6").Value  = 1.13
$ws.Range("AN6").Value = 15

# --- Row 7 updates ---
$ws.Range("G7").Value  = 4.33
$ws.Range("I7").Value  = 1.85
$ws.Range("J7").Value  = 5
$ws.Range("L7").Value  = 2.6
$ws.Range("O7").Value  = 1.4
$ws.Range("P7").Value  = 2.75
$ws.Range("U7").Value  = 3.25
$ws.Range("V7").Value  = 1.33
$ws.Range("Y7").Value  = 1.5
$ws.Range("Z7").Value  = 2.5
$ws.Range("AC7").Value = 10
$ws.Range("AD7").Value = 21
$ws.Range("AF7").Value = 51
$ws.Range("AH7").Value = 51
$ws.Range("AP7").Value = 15
$ws.Range("AQ7").Value = 17
$ws.Range("AS7").Value = 1000
